$d = $word.ActiveDocument

# Locate the paragraph that reads "Mammaw's Prune Cake" and remove the
# empty paragraphs immediately before and after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Mammaw.s Prune Cake") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $nextPara = $target.Next()
    if ($nextPara -ne $null -and $nextPara.Range.Text.Trim() -eq "") {
        $nextPara.Range.Delete()
    }

    $prevPara = $target.Previous()
    if ($prevPara -ne $null -and $prevPara.Range.Text.Trim() -eq "") {
        $prevPara.Range.Delete()
    }
}
